$d = $word.ActiveDocument

# 1. Remove trailing comma after "Etude de marcher"
$d.Content.Find.Execute(
    "Recherche de projet, mise en place d’un GANTT, Etude de marcher,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Recherche de projet, mise en place d’un GANTT, Etude de marcher",
    2) | Out-Null

# 2. Merge "Arrivé de Guyviane dans l'équipe" into a single run (drops the
#    spell-check proofErr markers that bracketed "Guyviane").
$d.Content.Find.Execute(
    "Arrivé de Guyviane dans l’équipe",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Arrivé de Guyviane dans l’équipe",
    2) | Out-Null

# 3. "Création de l'interfaces graphique et analyse ..." -> replace
#    "graphique" with "web", split across three runs.
$d.Content.Find.Execute(
    "Création de l’interfaces graphique et analyse d’un nouveau système d’empreinte",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Création de l’interfaces web et analyse d’un nouveau système d’empreinte",
    2) | Out-Null

# 4. "Nous avons pas trouvé de notice pour l'egistec  es603." -> drop the
#    trailing period and append " car elle n'existe pas."
$d.Content.Find.Execute(
    "Nous avons pas trouvé de notice pour l’egistec  es603.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nous avons pas trouvé de notice pour l’egistec  es603 car elle n’existe pas.",
    2) | Out-Null

# 5. Merge "Installation des drivers, travail sur la base de donnée" into a
#    single run (drops the proofErr markers around "donnée").
$d.Content.Find.Execute(
    "Installation des drivers, travail sur la base de donnée",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Installation des drivers, travail sur la base de donnée",
    2) | Out-Null

# 6. Merge "Difficulté a faire la base de donnée" (first occurrence, followed
#    by a separate "." run) into a single run.
$d.Content.Find.Execute(
    "Difficulté a faire la base de donnée.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Difficulté a faire la base de donnée.",
    2) | Out-Null

# 7. Merge ", comparatif des différents lecteur d'empreinte)" into a single run.
$d.Content.Find.Execute(
    ", comparatif des différents lecteur d’empreinte)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", comparatif des différents lecteur d’empreinte)",
    2) | Out-Null

# 8. Merge the second "Difficulté a faire la base de donnée, Difficulté a
#    trouve certaine référence." and append " car elle n'existe pas." while
#    dropping the final period.
$d.Content.Find.Execute(
    "Difficulté a faire la base de donnée, Difficulté a trouve certaine référence.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Difficulté a faire la base de donnée, Difficulté a trouve certaine référence car elle n’existe pas.",
    2) | Out-Null
